$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema3g"
$ws.Cells.Item(2,3).Value = "Nrp2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 11.42569233333333
$ws.Cells.Item(2,8).Value = 34.277077
$ws.Cells.Item(2,9).Value = 0.9287201125191051
$ws.Cells.Item(2,10).Value = 0.9504172085657416
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 44.13164066666667
$ws.Cells.Item(2,14).Value = 132.394922
$ws.Cells.Item(2,15).Value = 0.4415399811720331
$ws.Cells.Item(2,16).Value = 0.4562856844211927
$ws.Cells.Item(2,17).Value = 504.2345484225549
$ws.Cells.Item(2,18).Value = 4538.110935802994
$ws.Cells.Item(2,19).Value = 0.4100670609957742
$ws.Cells.Item(2,20).Value = 0.4336617664960989

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema3g"
$ws.Cells.Item(3,3).Value = "Nrp2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 11.42569233333333
$ws.Cells.Item(3,8).Value = 34.277077
$ws.Cells.Item(3,9).Value = 0.9287201125191051
$ws.Cells.Item(3,10).Value = 0.9504172085657416
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 14.93259333333333
$ws.Cells.Item(3,14).Value = 44.79778
$ws.Cells.Item(3,15).Value = 0.1494015830739255
$ws.Cells.Item(3,16).Value = 0.1543910098595022
$ws.Cells.Item(3,17).Value = 170.6152171654511
$ws.Cells.Item(3,18).Value = 1535.53695448906
$ws.Cells.Item(3,19).Value = 0.1387522550429485
$ws.Cells.Item(3,20).Value = 0.1467358726183139

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema3g"
$ws.Cells.Item(4,3).Value = "Nrp2"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 11.42569233333333
$ws.Cells.Item(4,8).Value = 34.277077
$ws.Cells.Item(4,9).Value = 0.9287201125191051
$ws.Cells.Item(4,10).Value = 0.9504172085657416
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 13.61024133333333
$ws.Cells.Item(4,14).Value = 40.830724
$ws.Cells.Item(4,15).Value = 0.1361713639304118
$ws.Cells.Item(4,16).Value = 0.1407189532975654
$ws.Cells.Item(4,17).Value = 155.5064300570831
$ws.Cells.Item(4,18).Value = 1399.557870513748
$ws.Cells.Item(4,19).Value = 0.126465084431332
$ws.Cells.Item(4,20).Value = 0.1337417147853651

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Sema3g"
$ws.Cells.Item(5,3).Value = "Nrp2"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 11.42569233333333
$ws.Cells.Item(5,8).Value = 34.277077
$ws.Cells.Item(5,9).Value = 0.9287201125191051
$ws.Cells.Item(5,10).Value = 0.9504172085657416
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 17.584752
$ws.Cells.Item(5,14).Value = 52.754256
$ws.Cells.Item(5,15).Value = 0.1759366057935712
$ws.Cells.Item(5,16).Value = 0.1818121982434553
$ws.Cells.Item(5,17).Value = 200.9179661099679
$ws.Cells.Item(5,18).Value = 1808.261694989712
$ws.Cells.Item(5,19).Value = 0.1633958643288349
$ws.Cells.Item(5,20).Value = 0.1727974419377461

$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Sema3g"
$ws.Cells.Item(6,3).Value = "Nrp2"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 11.42569233333333
$ws.Cells.Item(6,8).Value = 34.277077
$ws.Cells.Item(6,9).Value = 0.9287201125191051
$ws.Cells.Item(6,10).Value = 0.9504172085657416
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 9.6901375
$ws.Cells.Item(6,14).Value = 19.380275
$ws.Cells.Item(6,15).Value = 0.09695046603005844
$ws.Cells.Item(6,16).Value = 0.06679215417828435
$ws.Cells.Item(6,17).Value = 110.7165297426958
$ws.Cells.Item(6,18).Value = 664.299178456175
$ws.Cells.Item(6,19).Value = 0.09003984772021556
$ws.Cells.Item(6,20).Value = 0.06348041272821765

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema3g"
$ws.Cells.Item(7,3).Value = "Nrp2"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.034359
$ws.Cells.Item(7,8).Value = 0.103077
$ws.Cells.Item(7,9).Value = 0.00279281932465046
$ws.Cells.Item(7,10).Value = 0.00285806618246156
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 44.13164066666667
$ws.Cells.Item(7,14).Value = 132.394922
$ws.Cells.Item(7,15).Value = 0.4415399811720331
$ws.Cells.Item(7,16).Value = 0.4562856844211927
$ws.Cells.Item(7,17).Value = 1.516319041666
$ws.Cells.Item(7,18).Value = 13.646871374994
$ws.Cells.Item(7,19).Value = 0.001233141392023054
$ws.Cells.Item(7,20).Value = 0.001304094684185539

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Sema3g"
$ws.Cells.Item(8,3).Value = "Nrp2"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.034359
$ws.Cells.Item(8,8).Value = 0.103077
$ws.Cells.Item(8,9).Value = 0.00279281932465046
$ws.Cells.Item(8,10).Value = 0.00285806618246156
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 14.93259333333333
$ws.Cells.Item(8,14).Value = 44.79778
$ws.Cells.Item(8,15).Value = 0.1494015830739255
$ws.Cells.Item(8,16).Value = 0.1543910098595022
$ws.Cells.Item(8,17).Value = 0.5130689743400001
$ws.Cells.Item(8,18).Value = 4.61762076906
$ws.Cells.Item(8,19).Value = 0.0004172516283422301
$ws.Cells.Item(8,20).Value = 0.0004412597241555325

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Sema3g"
$ws.Cells.Item(9,3).Value = "Nrp2"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.034359
$ws.Cells.Item(9,8).Value = 0.103077
$ws.Cells.Item(9,9).Value = 0.00279281932465046
$ws.Cells.Item(9,10).Value = 0.00285806618246156
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 13.61024133333333
$ws.Cells.Item(9,14).Value = 40.830724
$ws.Cells.Item(9,15).Value = 0.1361713639304118
$ws.Cells.Item(9,16).Value = 0.1407189532975654
$ws.Cells.Item(9,17).Value = 0.467634281972
$ws.Cells.Item(9,18).Value = 4.208708537748
$ws.Cells.Item(9,19).Value = 0.0003803020166488646
$ws.Cells.Item(9,20).Value = 0.0004021840816511594

$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Sema3g"
$ws.Cells.Item(10,3).Value = "Nrp2"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.034359
$ws.Cells.Item(10,8).Value = 0.103077
$ws.Cells.Item(10,9).Value = 0.00279281932465046
$ws.Cells.Item(10,10).Value = 0.00285806618246156
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 17.584752
$ws.Cells.Item(10,14).Value = 52.754256
$ws.Cells.Item(10,15).Value = 0.1759366057935712
$ws.Cells.Item(10,16).Value = 0.1818121982434553
$ws.Cells.Item(10,17).Value = 0.604194493968
$ws.Cells.Item(10,18).Value = 5.437750445712
$ws.Cells.Item(10,19).Value = 0.0004913591525736958
$ws.Cells.Item(10,20).Value = 0.0005196312953586168

$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Sema3g"
$ws.Cells.Item(11,3).Value = "Nrp2"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.034359
$ws.Cells.Item(11,8).Value = 0.103077
$ws.Cells.Item(11,9).Value = 0.00279281932465046
$ws.Cells.Item(11,10).Value = 0.00285806618246156
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 9.6901375
$ws.Cells.Item(11,14).Value = 19.380275
$ws.Cells.Item(11,15).Value = 0.09695046603005844
$ws.Cells.Item(11,16).Value = 0.06679215417828435
$ws.Cells.Item(11,17).Value = 0.3329434343625
$ws.Cells.Item(11,18).Value = 1.997660606175
$ws.Cells.Item(11,19).Value = 0.0002707651350626152
$ws.Cells.Item(11,20).Value = 0.0001908963971107131

$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Sema3g"
$ws.Cells.Item(12,3).Value = "Nrp2"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.8425705
$ws.Cells.Item(12,8).Value = 1.685141
$ws.Cells.Item(12,9).Value = 0.06848706815624436
$ws.Cells.Item(12,10).Value = 0.04672472525179678
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 44.13164066666667
$ws.Cells.Item(12,14).Value = 132.394922
$ws.Cells.Item(12,15).Value = 0.4415399811720331
$ws.Cells.Item(12,16).Value = 0.4562856844211927
$ws.Cells.Item(12,17).Value = 37.18401854233367
$ws.Cells.Item(12,18).Value = 223.104111254002
$ws.Cells.Item(12,19).Value = 0.03023977878423589
$ws.Cells.Item(12,20).Value = 0.02131982324090828

$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Sema3g"
$ws.Cells.Item(13,3).Value = "Nrp2"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.8425705
$ws.Cells.Item(13,8).Value = 1.685141
$ws.Cells.Item(13,9).Value = 0.06848706815624436
$ws.Cells.Item(13,10).Value = 0.04672472525179678
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 14.93259333333333
$ws.Cells.Item(13,14).Value = 44.79778
$ws.Cells.Item(13,15).Value = 0.1494015830739255
$ws.Cells.Item(13,16).Value = 0.1543910098595022
$ws.Cells.Item(13,17).Value = 12.58176263116333
$ws.Cells.Item(13,18).Value = 75.49057578698
$ws.Cells.Item(13,19).Value = 0.01023207640263474
$ws.Cells.Item(13,20).Value = 0.007213877517032686

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Sema3g"
$ws.Cells.Item(14,3).Value = "Nrp2"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.8425705
$ws.Cells.Item(14,8).Value = 1.685141
$ws.Cells.Item(14,9).Value = 0.06848706815624436
$ws.Cells.Item(14,10).Value = 0.04672472525179678
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 13.61024133333333
$ws.Cells.Item(14,14).Value = 40.830724
$ws.Cells.Item(14,15).Value = 0.1361713639304118
$ws.Cells.Item(14,16).Value = 0.1407189532975654
$ws.Cells.Item(14,17).Value = 11.46758784534733
$ws.Cells.Item(14,18).Value = 68.80552707208399
$ws.Cells.Item(14,19).Value = 0.009325977482430865
$ws.Cells.Item(14,20).Value = 0.006575054430549167

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Sema3g"
$ws.Cells.Item(15,3).Value = "Nrp2"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.8425705
$ws.Cells.Item(15,8).Value = 1.685141
$ws.Cells.Item(15,9).Value = 0.06848706815624436
$ws.Cells.Item(15,10).Value = 0.04672472525179678
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 17.584752
$ws.Cells.Item(15,14).Value = 52.754256
$ws.Cells.Item(15,15).Value = 0.1759366057935712
$ws.Cells.Item(15,16).Value = 0.1818121982434553
$ws.Cells.Item(15,17).Value = 14.816393285016
$ws.Cells.Item(15,18).Value = 88.898359710096
$ws.Cells.Item(15,19).Value = 0.01204938231216261
$ws.Cells.Item(15,20).Value = 0.00849512501035066

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Sema3g"
$ws.Cells.Item(16,3).Value = "Nrp2"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.8425705
$ws.Cells.Item(16,8).Value = 1.685141
$ws.Cells.Item(16,9).Value = 0.06848706815624436
$ws.Cells.Item(16,10).Value = 0.04672472525179678
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 9.6901375
$ws.Cells.Item(16,14).Value = 19.380275
$ws.Cells.Item(16,15).Value = 0.09695046603005844
$ws.Cells.Item(16,16).Value = 0.06679215417828435
$ws.Cells.Item(16,17).Value = 8.16462399844375
$ws.Cells.Item(16,18).Value = 32.658495993775
$ws.Cells.Item(16,19).Value = 0.006639853174780266
$ws.Cells.Item(16,20).Value = 0.003120845052955987

